# Weekly data refresh: insert a new, most-recent observation at the top
# of the Apio data block (row 99), pushing all existing rows down by one.
# Excel's Rows(...).Insert() shifts the existing data down and extends
# the used range automatically (dimension A1:R188 -> A1:R189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(99).Insert()

$ws.Cells.Item(99, 1).Value = 10
$ws.Cells.Item(99, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(99, 3).Value = "La Araucanía"
$ws.Cells.Item(99, 4).Value = 44484
$ws.Cells.Item(99, 5).Value = 9
$ws.Cells.Item(99, 6).Value = 100112017
$ws.Cells.Item(99, 7).Value = "Apio"
$ws.Cells.Item(99, 8).Value = "Americana (o)"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 60
$ws.Cells.Item(99, 11).Value = 9000
$ws.Cells.Item(99, 12).Value = 9000
$ws.Cells.Item(99, 13).Value = 9000
$ws.Cells.Item(99, 14).Value = "$/docena de matas"
$ws.Cells.Item(99, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(99, 16).Value = 1500
$ws.Cells.Item(99, 17).Value = 6
$ws.Cells.Item(99, 18).Value = "Hortaliza"
